# Generate Report for Handoff
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect that a new handoff has occurred (status "Ready for handoff",
# a new handoff file named after commit 63290e5768f688058c7b37413b0a5c26c308f864,
# and a new handoff datetime).

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, $cellAddr, $newText) {
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address(0, 0)
        if ($addr -eq $cellAddr) {
            $hl.TextToDisplay = $newText
        }
    }
}

# --- Overview sheet: row 3 is the "b.md" file ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-26-20 08:26:53"

# --- zh-cn sheet: row 3 is the "b.md" file ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-20 08:26:50"
Set-HyperlinkDisplay $zh "D3" "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# --- de-de sheet: row 3 is the "b.md" file ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("E3").Value = "2016-03-20 08:26:53"
Set-HyperlinkDisplay $de "D3" "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
